$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.776.88"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.863.63"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  -4.97%  "
$ws.Range("D6").Value = "241.27"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "0.3092"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "24.53"
$ws.Range("E9").Value = "  -4.27%  "
$ws.Range("D10").Value = "0.07046"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").Value = "0.08404"
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("D12").Value = "0.7466"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "1.874.37"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "5.315"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "92.08"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "29.789.77"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "6.006"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("D19").Value = "239.64"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "2.137.49"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "7.875"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").Value = "0.1556"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "9.265"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "162.29"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").Value = "1.992"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "1.489"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").Value = "1.525"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "4.445"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "4.126"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "0.7421"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "2.700"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "2.756"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "0.4418"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "1.102.21"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "5.981"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "71.72"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").Value = "0.8604"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "7.698"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "101.83"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "1.827"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").Value = "2.989"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "2.033.14"
$ws.Range("E51").Value = "  +0.94%  "
